$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add header cells I1 ("I0") and J1 ("IF"), matching the style of the
# existing header cells (bold, bordered, centered).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I and J, rows 2-15.
$values = @{
  2  = @(8, 8)
  3  = @(8, 8)
  4  = @(7, 7)
  5  = @(5, 5)
  6  = @(6, 6)
  7  = @(7, 7)
  8  = @(5, 5)
  9  = @(7, 7)
  10 = @(9, 9)
  11 = @(4, 4)
  12 = @(6, 7)
  13 = @(9, 9)
  14 = @(4, 4)
  15 = @(7, 7)
}

foreach ($row in $values.Keys) {
  $pair = $values[$row]
  $ws.Cells.Item($row, 9).Value = $pair[0]
  $ws.Cells.Item($row, 10).Value = $pair[1]
}
